$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1467.5555
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1467.5555
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4402.666499999999
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -4738.666499999999
$ws.Range("H62").Value = 4499
$ws.Range("I62").Value = 4448.9
$ws.Range("K62").Value = 4448.9
$ws.Range("M62").Value = -3824.9
$ws.Range("H65").Value = 4499
$ws.Range("I65").Value = 4448.9
$ws.Range("K65").Value = 22244.5
$ws.Range("M65").Value = -19124.5
$ws.Range("H86").Value = 2358.6667
$ws.Range("I86").Value = 1545
$ws.Range("J86").Value = 4799.6665
$ws.Range("K86").Value = 1545
$ws.Range("L86").Value = 4799.6665
$ws.Range("M86").Value = -422
$ws.Range("N86").Value = -7045.6665
$ws.Range("H89").Value = 2358.6667
$ws.Range("I89").Value = 1545
$ws.Range("J89").Value = 4799.6665
$ws.Range("K89").Value = 7725
$ws.Range("L89").Value = 23998.3325
$ws.Range("M89").Value = -2109
$ws.Range("N89").Value = -35230.3325
$ws.Range("H132").Value = 1602.2084
$ws.Range("I132").Value = 1596.381
$ws.Range("K132").Value = 4789.143
$ws.Range("M132").Value = -2259.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3293.4285
$ws.Range("J2").Value = 4899.3335
$ws.Range("L2").Value = 4899.3335
$ws.Range("N2").Value = -5125.3335
$ws.Range("H61").Value = 2401.6
$ws.Range("I61").Value = 2401.6
$ws.Range("K61").Value = 2401.6
$ws.Range("M61").Value = -2189.6
$ws.Range("H74").Value = 1438.6666
$ws.Range("I74").Value = 1114.4706
$ws.Range("K74").Value = 1114.4706
$ws.Range("M74").Value = -240.4706000000001
$ws.Range("H77").Value = 1438.6666
$ws.Range("I77").Value = 1114.4706
$ws.Range("K77").Value = 5572.353000000001
$ws.Range("M77").Value = -1204.353000000001
$ws.Range("H97").Value = 3217.6667
$ws.Range("J97").Value = 3703.3333
$ws.Range("L97").Value = 3703.3333
$ws.Range("N97").Value = -4695.3333
$ws.Range("H116").Value = 3293.4285
$ws.Range("J116").Value = 4899.3335
$ws.Range("L116").Value = 4899.3335
$ws.Range("N116").Value = -9487.333500000001
$ws.Range("H136").Value = 2401.6
$ws.Range("I136").Value = 2401.6
$ws.Range("K136").Value = 7204.799999999999
$ws.Range("M136").Value = -4654.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3293.4285
$ws.Range("J3").Value = 4899.3335
$ws.Range("L3").Value = 4899.3335
$ws.Range("N3").Value = -5127.3335
$ws.Range("H82").Value = 25085.666
$ws.Range("I82").Value = 25085.666
$ws.Range("K82").Value = 25085.666
$ws.Range("M82").Value = -24702.666
$ws.Range("H85").Value = 25085.666
$ws.Range("I85").Value = 25085.666
$ws.Range("K85").Value = 25085.666
$ws.Range("M85").Value = -23759.666
$ws.Range("H94").Value = 1160.2858
$ws.Range("J94").Value = 425
$ws.Range("L94").Value = 425
$ws.Range("N94").Value = -1327
$ws.Range("H99").Value = 45224.74
$ws.Range("I99").Value = 57187.168
$ws.Range("K99").Value = 57187.168
$ws.Range("M99").Value = -55689.168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 239
$ws.Range("I22").Value = 239
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 239
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 111
$ws.Range("N22").ClearContents()
$ws.Range("H38").Value = 3036.3333
$ws.Range("I38").Value = 3036.3333
$ws.Range("K38").Value = 3036.3333
$ws.Range("M38").Value = -2659.3333
$ws.Range("H46").Value = 3036.3333
$ws.Range("I46").Value = 3036.3333
$ws.Range("K46").Value = 3036.3333
$ws.Range("M46").Value = -2825.3333
$ws.Range("H99").Value = 11806.407
$ws.Range("I99").Value = 7255.0713
$ws.Range("J99").Value = 16707.846
$ws.Range("K99").Value = 7255.0713
$ws.Range("L99").Value = 16707.846
$ws.Range("M99").Value = -5757.0713
$ws.Range("N99").Value = -19703.846
$ws.Range("H126").Value = 11806.407
$ws.Range("I126").Value = 7255.0713
$ws.Range("J126").Value = 16707.846
$ws.Range("K126").Value = 21765.2139
$ws.Range("L126").Value = 50123.538
$ws.Range("M126").Value = -19295.2139
$ws.Range("N126").Value = -55063.538
$ws.Range("H132").Value = 4097.846
$ws.Range("I132").Value = 3025.0908
$ws.Range("K132").Value = 9075.2724
$ws.Range("M132").Value = -6545.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1072240.9
$ws.Range("I4").Value = 1372242.9
$ws.Range("J4").Value = 805.1429000000001
$ws.Range("K4").Value = 4116728.7
$ws.Range("L4").Value = 2415.4287
$ws.Range("M4").Value = -4116616.7
$ws.Range("N4").Value = -2639.4287
$ws.Range("H7").Value = 25000082
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H92").Value = 381.57144
$ws.Range("I92").Value = 390.33334
$ws.Range("J92").Value = 375
$ws.Range("K92").Value = 1171.00002
$ws.Range("L92").Value = 1125
$ws.Range("M92").Value = 76.99998000000005
$ws.Range("N92").Value = -3621
$ws.Range("H105").Value = 4000
$ws.Range("J105").Value = 4000
$ws.Range("L105").Value = 12000
$ws.Range("N105").Value = -17242
$ws.Range("H106").Value = 12999.8
$ws.Range("J106").Value = 12999.8
$ws.Range("L106").Value = 38999.39999999999
$ws.Range("N106").Value = -40891.39999999999
$ws.Range("H134").Value = 2027.4
$ws.Range("I134").Value = 2027.4
$ws.Range("K134").Value = 6082.200000000001
$ws.Range("M134").Value = -1012.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6880.067
$ws.Range("I70").Value = 5905.875
$ws.Range("K70").Value = 5905.875
$ws.Range("M70").Value = -5635.875
$ws.Range("H73").Value = 6880.067
$ws.Range("I73").Value = 5905.875
$ws.Range("K73").Value = 5905.875
$ws.Range("M73").Value = -4969.875
$ws.Range("H113").Value = 2969.5
$ws.Range("J113").Value = 2939.75
$ws.Range("L113").Value = 2939.75
$ws.Range("N113").Value = -7279.75
$ws.Range("H132").Value = 2760.122
$ws.Range("I132").Value = 2649
$ws.Range("J132").Value = 2847.087
$ws.Range("K132").Value = 7947
$ws.Range("L132").Value = 8541.261
$ws.Range("M132").Value = -5417
$ws.Range("N132").Value = -13601.261

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 2557750
$ws.Range("J43").Value = 2557750
$ws.Range("L43").Value = 2557750
$ws.Range("N43").Value = -2558136

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1901.25
$ws.Range("I100").Value = 1820.1666
$ws.Range("K100").Value = 3640.3332
$ws.Range("M100").Value = -3099.3332
$ws.Range("H132").Value = 42140.914
$ws.Range("I132").Value = 63874.6
$ws.Range("K132").Value = 191623.8
$ws.Range("M132").Value = -189093.8
$ws.Range("H136").Value = 1604.68
$ws.Range("I136").Value = 1551.7273
$ws.Range("K136").Value = 4655.1819
$ws.Range("M136").Value = -2105.1819
